# Insert a new data row at row 271 (pushing existing rows 271..364 down to 272..365)
# and populate it with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("271:271").Insert()

$ws.Range("A271").Value2 = 8
$ws.Range("B271").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C271").Value2 = "Coquimbo"
$ws.Range("D271").Value2 = 44468
$ws.Range("E271").Value2 = 4
$ws.Range("F271").Value2 = 100112024
$ws.Range("G271").Value2 = "Choclo"
$ws.Range("H271").Value2 = "Dulce o Americano"
$ws.Range("I271").Value2 = "Primera"
$ws.Range("J271").Value2 = 560
$ws.Range("K271").Value2 = 37000
$ws.Range("L271").Value2 = 38000
$ws.Range("M271").Value2 = 37500
$ws.Range("N271").Value2 = "$/malla 70 unidades"
$ws.Range("O271").Value2 = "Región de Arica y Parinacota"
$ws.Range("P271").Value2 = 536
$ws.Range("Q271").Value2 = 70
$ws.Range("R271").Value2 = "Hortaliza"
